$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 1047; this shifts the existing rows
# 1047..1148 down to 1048..1149 and expands the used range/dimension
# to A1:T1149 automatically.
$ws.Rows(1047).Insert()

# Populate the newly inserted row 1047 with the new data record.
$ws.Range("A1047").Value2 = 3
$ws.Range("B1047").Value2 = "Femacal de La Calera"
$ws.Range("C1047").Value2 = "Coquimbo"
$ws.Range("D1047").Value2 = 44946
$ws.Range("E1047").Value2 = 5
$ws.Range("F1047").Value2 = "Fruta"
$ws.Range("G1047").Value2 = 100102
$ws.Range("H1047").Value2 = "Cítricos"
$ws.Range("I1047").Value2 = 100102005
$ws.Range("J1047").Value2 = "Naranja"
$ws.Range("K1047").Value2 = "Valencia"
$ws.Range("L1047").Value2 = "Primera"
$ws.Range("M1047").Value2 = 133
$ws.Range("N1047").Value2 = 8500
$ws.Range("O1047").Value2 = 9000
$ws.Range("P1047").Value2 = 8756
$ws.Range("Q1047").Value2 = "$/malla 13 kilos"
$ws.Range("R1047").Value2 = "Provincia de Quillota"
$ws.Range("S1047").Value2 = 674
$ws.Range("T1047").Value2 = 13
